$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D holds price figures that are stored as TEXT in this sheet
# (not numbers) even when they look like plain decimals, e.g. "206.44".
# Assigning a bare numeric-looking string via .Value would make Excel
# auto-convert it to a real number (and drop formatting like trailing
# zeros), so those assignments below use a leading apostrophe - exactly
# like typing '206.36 into a cell in the Excel UI - to force text storage
# while still landing the correct visible characters in the cell.

# --- Simple value updates (price / volume columns), rows 2-43 ---
$ws.Range("D2").Value = "27.160.49"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "1.558.80"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'206.36"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'0.490"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'22.13"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "'0.0860"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "1.781.12"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "1.558.18"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "'62.84"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "27.164.10"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "'214.09"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").Value = "'7.22"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "'9.33"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "'151.89"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").Value = "'6.59"
$ws.Range("E26").Value = "  -3.37%  "
$ws.Range("D27").Value = "'14.87"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "1.378.98"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "'2.92"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'0.948"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "'0.809"
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "'0.985"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").Value = "'1.80"
$ws.Range("E43").Value = "  +3.74%  "

# --- Rows 44 and 45 swap places (MXToken moves up to rank 44, Aave drops
#     to rank 45), each bringing its own refreshed price / volume figures.
#     Rank numbers in column A (42 / 43) stay as-is. ---
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.16"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'63.16"
$ws.Range("E45").Value = "  -1.72%  "

# --- Remaining value updates, rows 46-51 ---
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "1.693.18"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("D48").Value = "'85.33"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").Value = "0.0₇0983"
$ws.Range("E49").Value = "  -2.90%  "
$ws.Range("D50").Value = "'0.0492"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("E51").Value = "  +0.20%  "
